# Update Sheets via scheduled runner
# Applies updated market price / profit figures to the Leve profit tables
# across all class sheets, per the latest data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4198.7
$ws.Range("I51").Value = 2183.75
$ws.Range("J51").Value = 5542
$ws.Range("K51").Value = 2183.75
$ws.Range("L51").Value = 5542
$ws.Range("M51").Value = -1699.75
$ws.Range("N51").Value = -6510
$ws.Range("H62").Value = 3200.5557
$ws.Range("I62").Value = 3286.4285
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 3286.4285
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -2662.4285
$ws.Range("N62").Value = -4148
$ws.Range("H65").Value = 3200.5557
$ws.Range("I65").Value = 3286.4285
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 16432.1425
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -13312.1425
$ws.Range("N65").Value = -20740
$ws.Range("H87").Value = 33656
$ws.Range("J87").Value = 33656
$ws.Range("L87").Value = 33656
$ws.Range("N87").Value = -36152
$ws.Range("H90").Value = 33656
$ws.Range("J90").Value = 33656
$ws.Range("L90").Value = 100968
$ws.Range("N90").Value = -113448
$ws.Range("H111").Value = 2058.6
$ws.Range("I111").Value = 1617.091
$ws.Range("J111").Value = 2598.2222
$ws.Range("K111").Value = 4851.272999999999
$ws.Range("L111").Value = 7794.6666
$ws.Range("M111").Value = -1784.272999999999
$ws.Range("N111").Value = -13928.6666
$ws.Range("H113").Value = 3772.913
$ws.Range("I113").Value = 3637
$ws.Range("J113").Value = 3921.182
$ws.Range("K113").Value = 3637
$ws.Range("L113").Value = 3921.182
$ws.Range("M113").Value = -383
$ws.Range("N113").Value = -10429.182
$ws.Range("H129").Value = 682.25
$ws.Range("J129").Value = 998
$ws.Range("L129").Value = 2994
$ws.Range("N129").Value = -12994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1517.2307
$ws.Range("I45").Value = 1553
$ws.Range("J45").Value = 1460
$ws.Range("K45").Value = 1553
$ws.Range("L45").Value = 1460
$ws.Range("M45").Value = -1176
$ws.Range("N45").Value = -2214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 24500
$ws.Range("J53").Value = 24500
$ws.Range("L53").Value = 24500
$ws.Range("N53").Value = -25648
$ws.Range("H75").Value = 4394.5
$ws.Range("I75").Value = 1131.2
$ws.Range("J75").Value = 9833.333000000001
$ws.Range("K75").Value = 1131.2
$ws.Range("L75").Value = 9833.333000000001
$ws.Range("M75").Value = -195.2
$ws.Range("N75").Value = -11705.333
$ws.Range("H78").Value = 4394.5
$ws.Range("I78").Value = 1131.2
$ws.Range("J78").Value = 9833.333000000001
$ws.Range("K78").Value = 3393.6
$ws.Range("L78").Value = 29499.999
$ws.Range("M78").Value = 1286.4
$ws.Range("N78").Value = -38859.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 51271.89
$ws.Range("I31").Value = 43880.707
$ws.Range("J31").Value = 59718.953
$ws.Range("K31").Value = 43880.707
$ws.Range("L31").Value = 59718.953
$ws.Range("M31").Value = -43585.707
$ws.Range("N31").Value = -60308.953
$ws.Range("H34").Value = 51271.89
$ws.Range("I34").Value = 43880.707
$ws.Range("J34").Value = 59718.953
$ws.Range("K34").Value = 43880.707
$ws.Range("L34").Value = 59718.953
$ws.Range("M34").Value = -43678.707
$ws.Range("N34").Value = -60122.953
$ws.Range("H58").Value = 27779864
$ws.Range("I58").Value = 41668496
$ws.Range("J58").Value = 2598.75
$ws.Range("K58").Value = 41668496
$ws.Range("L58").Value = 2598.75
$ws.Range("M58").Value = -41668293
$ws.Range("N58").Value = -3004.75
$ws.Range("H92").Value = 20280.2
$ws.Range("J92").Value = 20280.2
$ws.Range("L92").Value = 20280.2
$ws.Range("N92").Value = -25272.2
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800
$ws.Range("H136").Value = 27779864
$ws.Range("I136").Value = 41668496
$ws.Range("J136").Value = 2598.75
$ws.Range("K136").Value = 125005488
$ws.Range("L136").Value = 7796.25
$ws.Range("M136").Value = -125002938
$ws.Range("N136").Value = -12896.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 839.0125
$ws.Range("J131").Value = 923.8677
$ws.Range("L131").Value = 2771.6031
$ws.Range("N131").Value = -12851.6031
$ws.Range("H134").Value = 7559.6875
$ws.Range("I134").Value = 4777.5
$ws.Range("J134").Value = 7957.143
$ws.Range("K134").Value = 14332.5
$ws.Range("L134").Value = 23871.429
$ws.Range("M134").Value = -9262.5
$ws.Range("N134").Value = -34011.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 22000
$ws.Range("I93").Value = 22000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 22000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -20128
$ws.Range("N93").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 2146.4614
$ws.Range("I122").Value = 1926.5294
$ws.Range("J122").Value = 2561.889
$ws.Range("K122").Value = 5779.5882
$ws.Range("L122").Value = 7685.667
$ws.Range("M122").Value = -3329.5882
$ws.Range("N122").Value = -12585.667
$ws.Range("H132").Value = 69246.164
$ws.Range("I132").Value = 41854.36
$ws.Range("J132").Value = 206205.2
$ws.Range("K132").Value = 125563.08
$ws.Range("L132").Value = 618615.6000000001
$ws.Range("M132").Value = -123033.08
$ws.Range("N132").Value = -623675.6000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3765.641
$ws.Range("I122").Value = 4004.2354
$ws.Range("J122").Value = 3581.2727
$ws.Range("K122").Value = 12012.7062
$ws.Range("L122").Value = 10743.8181
$ws.Range("M122").Value = -9562.706200000001
$ws.Range("N122").Value = -15643.8181
$ws.Range("H132").Value = 25243.37
$ws.Range("I132").Value = 2817.3572
$ws.Range("K132").Value = 8452.071599999999
$ws.Range("M132").Value = -5922.071599999999
$ws.Range("H136").Value = 46400
$ws.Range("I136").Value = 28165.879
$ws.Range("J136").Value = 170999.83
$ws.Range("K136").Value = 84497.637
$ws.Range("L136").Value = 512999.49
$ws.Range("M136").Value = -81947.637
$ws.Range("N136").Value = -518099.49
$ws.Range("H139").Value = 38655.266
$ws.Range("J139").Value = 38655.266
$ws.Range("L139").Value = 38655.266
$ws.Range("N139").Value = -48935.266

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1509
$ws.Range("I122").Value = 1190.9231
$ws.Range("J122").Value = 2198.1667
$ws.Range("K122").Value = 3572.7693
$ws.Range("L122").Value = 6594.500100000001
$ws.Range("M122").Value = -1122.7693
$ws.Range("N122").Value = -11494.5001
$ws.Range("H132").Value = 40093.73
$ws.Range("I132").Value = 39599.08
$ws.Range("J132").Value = 40588.383
$ws.Range("K132").Value = 118797.24
$ws.Range("L132").Value = 121765.149
$ws.Range("M132").Value = -116267.24
$ws.Range("N132").Value = -126825.149
$ws.Range("H136").Value = 39251.203
$ws.Range("I136").Value = 24982.262
$ws.Range("J136").Value = 89192.5
$ws.Range("K136").Value = 74946.78599999999
$ws.Range("L136").Value = 267577.5
$ws.Range("M136").Value = -72396.78599999999
$ws.Range("N136").Value = -272677.5
